# Update cryptos list: prices (column D) and 1h volume % (column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.060.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "'2.468.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'583.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").Value = "'169.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "'2.468.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("D14").Value = "'25.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "'66.622.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("E17").Value = "  -4.72%  "
$ws.Range("D18").Value = "'2.456.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("D19").Value = "'11.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.67%  "
$ws.Range("E20").Value = "  -5.21%  "
$ws.Range("D21").Value = "'352.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("D22").Value = "'4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'69.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").Value = "'4.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.39%  "
$ws.Range("D26").Value = "'1.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.30%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.35%  "
$ws.Range("E28").Value = "  -64.58%  "
$ws.Range("D29").Value = "'2.567.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("D30").Value = "'517.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").Value = "'0.0₃0904"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.05%  "
$ws.Range("D32").Value = "'7.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.78%  "
$ws.Range("D33").Value = "'1.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("E34").Value = "  -6.75%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  -8.47%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D39").Value = "'18.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("E40").Value = "  -6.63%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -6.53%  "
$ws.Range("E43").Value = "  -7.09%  "
$ws.Range("E44").Value = "  -7.48%  "
$ws.Range("D45").Value = "'2.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.51%  "
$ws.Range("D46").Value = "'38.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("D47").Value = "'141.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.72%  "
$ws.Range("D48").Value = "'3.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.08%  "
$ws.Range("D49").Value = "'0.515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.23%  "
$ws.Range("D50").Value = "'0.0₆0255"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -13.29%  "
$ws.Range("E51").Value = "  -7.64%  "
